$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column B (casos) per semana epidemiologica 43 de 2025
$ws.Range("B20").Value = 400
$ws.Range("B21").Value = 303
$ws.Range("B25").Value = 261
$ws.Range("B26").Value = 351
$ws.Range("B27").Value = 273
$ws.Range("B28").Value = 237
$ws.Range("B30").Value = 341
$ws.Range("B31").Value = 380
$ws.Range("B33").Value = 258
$ws.Range("B37").Value = 447
$ws.Range("B44").Value = 441

# Add new row 45 (semana 44, casos 1)
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 1
